$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170 (shifts existing rows 170..259 down to 171..260)
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new record
$ws.Cells.Item(170, 1).Value = 5
$ws.Cells.Item(170, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(170, 3).Value = "Maule"
$ws.Cells.Item(170, 4).Value = 44572
$ws.Cells.Item(170, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(170, 5).Value = 7
$ws.Cells.Item(170, 6).Value = 100114013
$ws.Cells.Item(170, 7).Value = "Zanahoria"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 300
$ws.Cells.Item(170, 11).Value = 7000
$ws.Cells.Item(170, 12).Value = 7000
$ws.Cells.Item(170, 13).Value = 7000
$ws.Cells.Item(170, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(170, 15).Value = "Región de Ñuble"
$ws.Cells.Item(170, 16).Value = 350
$ws.Cells.Item(170, 17).Value = 20
$ws.Cells.Item(170, 18).Value = "Hortaliza"
